$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(3,1,66.85111099999999,200.553333,0.1215550702639512,0.1215550702639512,3,1,3.024839333333334,9.074518000000001,0.1801507982970389,0.1801507982970388,202.2138700298327,1819.924830268494,0.02189824294510346,0.02189824294510346),
    @(3,1,66.85111099999999,200.553333,0.1215550702639512,0.1215550702639512,3,1,1.626140333333333,4.878420999999999,0.09684827751501936,0.09684827751501934,108.7092879252436,978.3835913271928,0.01177239917828083,0.01177239917828082),
    @(3,1,66.85111099999999,200.553333,0.1215550702639512,0.1215550702639512,3,1,4.640628666666666,13.921886,0.2763825997921178,0.2763825997921177,310.2311821051152,2792.080638946038,0.03359570633746439,0.03359570633746438),
    @(3,1,66.85111099999999,200.553333,0.1215550702639512,0.1215550702639512,3,1,7.498988000000001,22.496964,0.4466183243958241,0.446618324395824,501.315679175668,4511.841112581012,0.05428872180310255,0.05428872180310255),
    @(3,1,99.86393233333332,299.591797,0.1815821327429069,0.1815821327429069,3,1,3.024839333333334,9.074518000000001,0.1801507982970389,0.1801507982970388,302.0723505032051,2718.651154528846,0.03271216617011356,0.03271216617011356),
    @(3,1,99.86393233333332,299.591797,0.1815821327429069,0.1815821327429069,3,1,1.626140333333333,4.878420999999999,0.09684827751501936,0.09684827751501934,162.3927682125041,1461.534913912537,0.01758591678365413,0.01758591678365413),
    @(3,1,99.86393233333332,299.591797,0.1815821327429069,0.1815821327429069,3,1,4.640628666666666,13.921886,0.2763825997921178,0.2763825997921177,463.4314271521268,4170.882844369142,0.05018614192328204,0.05018614192328203),
    @(3,1,99.86393233333332,299.591797,0.1815821327429069,0.1815821327429069,3,1,7.498988000000001,22.496964,0.4466183243958241,0.446618324395824,748.8784302004786,6739.905871804308,0.08109790786585719,0.08109790786585717),
    @(3,1,368.4456226666667,1105.336868,0.6699429954379058,0.6699429954379058,3,1,3.024839333333334,9.074518000000001,0.1801507982970389,0.1801507982970388,1114.488811636625,10030.39930472963,0.1206907654416482,0.1206907654416482),
    @(3,1,368.4456226666667,1105.336868,0.6699429954379058,0.6699429954379058,3,1,1.626140333333333,4.878420999999999,0.09684827751501936,0.09684827751501934,599.1442876583808,5392.298588925428,0.06488282514141365,0.06488282514141364),
    @(3,1,368.4456226666667,1105.336868,0.6699429954379058,0.6699429954379058,3,1,4.640628666666666,13.921886,0.2763825997921178,0.2763825997921177,1709.819318654783,15388.37386789305,0.1851605867916473,0.1851605867916473),
    @(3,1,368.4456226666667,1105.336868,0.6699429954379058,0.6699429954379058,3,1,7.498988000000001,22.496964,0.4466183243958241,0.446618324395824,2762.969303029862,24866.72372726876,0.2992088180631967,0.2992088180631967),
    @(3,1,14.80496566666667,44.414897,0.02691980155523597,0.02691980155523597,3,1,3.024839333333334,9.074518000000001,0.1801507982970389,0.1801507982970388,44.7826424771829,403.0437822946461,0.004849623740173629,0.004849623740173627),
    @(3,1,14.80496566666667,44.414897,0.02691980155523597,0.02691980155523597,3,1,1.626140333333333,4.878420999999999,0.09684827751501936,0.09684827751501934,24.07495180418189,216.674566237637,0.002607136411670743,0.002607136411670743),
    @(3,1,14.80496566666667,44.414897,0.02691980155523597,0.02691980155523597,3,1,4.640628666666666,13.921886,0.2763825997921178,0.2763825997921177,68.70434808174912,618.3391327357419,0.007440164739724013,0.00744016473972401),
    @(3,1,14.80496566666667,44.414897,0.02691980155523597,0.02691980155523597,3,1,7.498988000000001,22.496964,0.4466183243958241,0.446618324395824,111.0222598747453,999.2003388727081,0.01202287666366759,0.01202287666366758),
)

$startRow = 2
$startCol = 5  # column E

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($r, $startCol + $j).Value = $rowVals[$j]
    }
}